# Golem_Profits scheduled update: refresh market-board derived columns (H-N)
# for a batch of Leve rows across the crafting-class sheets.
$wb = $excel.ActiveWorkbook

# ALC!132 - Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 908.1539
$ws.Range("I132").Value = 908.1539
$ws.Range("K132").Value = 2724.4617
$ws.Range("M132").Value = -194.4616999999998

# ALC!138 - All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2550.0588
$ws.Range("J138").Value = 3349.3
$ws.Range("L138").Value = 10047.9
$ws.Range("N138").Value = -20327.9

# ARM!38 - Eyes on a Hard Body
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 725848.3
$ws.Range("I38").Value = 3234.5
$ws.Range("J38").Value = 1689333.4
$ws.Range("K38").Value = 3234.5
$ws.Range("L38").Value = 1689333.4
$ws.Range("M38").Value = -2767.5
$ws.Range("N38").Value = -1690267.4

# ARM!61 - Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788

# ARM!112 - Wrapped Knuckles
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 5387
$ws.Range("J112").Value = 5387
$ws.Range("L112").Value = 5387
$ws.Range("N112").Value = -8341

# ARM!136 - Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

# BSM!94 - High Steal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2248.4707
$ws.Range("I94").Value = 1838.6364
$ws.Range("K94").Value = 1838.6364
$ws.Range("M94").Value = -1387.6364

# BSM!134 - Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2345.182
$ws.Range("I134").Value = 2345.182
$ws.Range("K134").Value = 7035.545999999999
$ws.Range("M134").Value = -4500.545999999999

# CRP!48 - The Cold, Cold Ground
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# CRP!69 - Landing the Big One
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 24333.334
$ws.Range("I69").Value = 10500
$ws.Range("K69").Value = 10500
$ws.Range("M69").Value = -9751

# CRP!72 - Fishing for Profits (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 24333.334
$ws.Range("I72").Value = 10500
$ws.Range("K72").Value = 31500
$ws.Range("M72").Value = -27756

# CRP!132 - Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2300.818
$ws.Range("I132").Value = 2039.125
$ws.Range("K132").Value = 6117.375
$ws.Range("M132").Value = -3587.375

# CUL!2 - Pork Is a Salty Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 319.57144
$ws.Range("I2").Value = 34.5
$ws.Range("K2").Value = 207
$ws.Range("M2").Value = -94

# CUL!14 - Keep Your Powder Dry
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1961.8572
$ws.Range("I14").Value = 1961.8572
$ws.Range("K14").Value = 5885.571599999999
$ws.Range("M14").Value = -5712.571599999999

# CUL!50 - Moving Up in the World
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 497.5
$ws.Range("I50").Value = 497.5
$ws.Range("K50").Value = 1492.5
$ws.Range("M50").Value = -1011.5

# CUL!53 - Rolanberry Fields Forever
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 497.5
$ws.Range("I53").Value = 497.5
$ws.Range("K53").Value = 1492.5
$ws.Range("M53").Value = -1011.5

# CUL!60 - Drinking to Your Health
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 267.66666
$ws.Range("I60").Value = 174
$ws.Range("K60").Value = 522
$ws.Range("M60").Value = -271

# CUL!103 - West Meats East
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 238.4
$ws.Range("I103").Value = 112
$ws.Range("J103").Value = 270
$ws.Range("K103").Value = 336
$ws.Range("L103").Value = 810
$ws.Range("M103").Value = 543
$ws.Range("N103").Value = -2568

# GSM!26 - Perk of Fiction
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("N26").ClearContents()

# GSM!35 - Necklet of Champions
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 5000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -5596

# GSM!50 - Coral on My Mind
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").ClearContents()

# GSM!70 - Sky Is the Limit
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 4000
$ws.Range("M70").Value = -3730

# GSM!73 - Hulls of Broken Dreams (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 4000
$ws.Range("M73").Value = -3064

# GSM!103 - Ring in the New
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("N103").ClearContents()

# GSM!126 - Gold Rush Order
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3363
$ws.Range("J126").Value = 4999
$ws.Range("L126").Value = 14997
$ws.Range("N126").Value = -19937

# LTW!31 - Open to Attack
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 7333
$ws.Range("I31").Value = 4999.5
$ws.Range("J31").Value = 12000
$ws.Range("K31").Value = 4999.5
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = -4751.5
$ws.Range("N31").Value = -12496

# LTW!32 - Men Who Scare Up Goats
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5141.5
$ws.Range("I32").Value = 1169.8
$ws.Range("K32").Value = 1169.8
$ws.Range("M32").Value = -852.8

# LTW!33 - Just Rewards
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 2000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -2580

# LTW!46 - Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3255
$ws.Range("I46").Value = 1006.6667
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 1006.6667
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -818.6667
$ws.Range("N46").Value = -10376

# LTW!56 - Hold On Tight
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").ClearContents()

# LTW!68 - You Could Say It's a Moving Target
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -3498

# LTW!71 - They Call It Bloody Mary (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -17488

# LTW!74 - Overall, We Blend In
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 15000
$ws.Range("J74").Value = 15000
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -16996

# LTW!77 - Eviction Notice (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 15000
$ws.Range("J77").Value = 15000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -54984

# LTW!122 - Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 300
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -5200

# WVR!54 - No Country for Cold Men
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 49998
$ws.Range("J54").Value = 49998
$ws.Range("L54").Value = 49998
$ws.Range("N54").Value = -51038

# WVR!58 - Seeing It Through to the End
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20616

# WVR!126 - A Polished Purchase
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2190.8333
$ws.Range("I126").Value = 2190.8333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6572.499899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4102.499899999999
$ws.Range("N126").ClearContents()
